# Update the "想去人数" (number of interested people) and "最低票价" (lowest ticket price)
# values for several rows on both the "展览" and "全部类型" worksheets.
# Both sheets contain the same underlying data, so the same edits are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: lowest ticket price 49.9 -> 65
    $ws.Range("G2").Value = 65

    # F4: interested count 289 -> 290
    $ws.Range("F4").Value = 290

    # F5: interested count 837 -> 838
    $ws.Range("F5").Value = 838

    # F6: interested count 6 -> 7
    $ws.Range("F6").Value = 7

    # F8: interested count 7476 -> 7561
    $ws.Range("F8").Value = 7561

    # F10: interested count 77 -> 78
    $ws.Range("F10").Value = 78

    # F11: interested count 127 -> 128
    $ws.Range("F11").Value = 128

    # F12: interested count 102 -> 103
    $ws.Range("F12").Value = 103

    # F17: interested count 26 -> 27
    $ws.Range("F17").Value = 27

    # F18: interested count 248 -> 249
    $ws.Range("F18").Value = 249

    # F19: interested count 677 -> 679
    $ws.Range("F19").Value = 679
}
